$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date-style) from row 5's A and G cells into row 6's A and G cells
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate new row 6 with trade data
$ws.Range("A6").Value = 42649.654224537036
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 10010.959999999999
$ws.Range("D6").Value = 10015.469999999999
$ws.Range("E6").Value = 77.349997999999999
$ws.Range("F6").Value = 77.42
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 0.09
$ws.Range("I6").Value = $false

# Widen column E slightly to fit new content (target stored width 9.875;
# the COM ColumnWidth setter quantizes to 1/6-character pixel steps, so 9.0
# is the closest input that lands on the nearest achievable stored width)
$ws.Columns.Item(5).ColumnWidth = 9
